$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.557.98'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '1.841.12'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '259.17'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5230'
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3174'
$ws.Range('E8').Value = '  -2.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06788'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.72'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7805'
$ws.Range('E11').Value = '  +2.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07772'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '1.827.96'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.94'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.011'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.85'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007934'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '26.584.72'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = '2.065.32'
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.606'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.965'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.334'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.67'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.219'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.674'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.92'
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.178'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08728'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.071'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04889'
$ws.Range('E33').Value = '  +2.17%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.132'
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7212'
$ws.Range('E35').Value = '  +2.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.858'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.096'
$ws.Range('E37').Value = '  +1.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.224'
$ws.Range('E38').Value = '  +1.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01742'
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4815'
$ws.Range('E40').Value = '  -0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8980'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '110.18'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.922'
$ws.Range('E43').Value = '  -2.31%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.650'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4160'
$ws.Range('E46').Value = '  +1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.000'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1232'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05833'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.86'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8929'
$ws.Range('E51').Value = '  +0.92%  '
